# C5-PowerPoint.pptx edit: 2020-05-04
#
# 1) Slide 6's table switches to a new built-in table style.
# 2) The deck's theme colour palette is swapped from "Integral" to the
#    stock "Office" palette (the slide master's theme, reached through
#    any slide's ThemeColorScheme, is the only theme resource PowerPoint
#    exposes through automation - the notes-master's theme part is not
#    reachable via the object model and is intentionally left alone).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table
$tbl.ApplyStyle("{5586E2D8-5A8E-4B29-924F-3D868D60B7A3}")

# --- 2) Theme colours -------------------------------------------------
# RGB is stored/returned as a BGR-packed long (the classic VBA RGB()
# layout), so 0xRRGGBB in the OOXML <a:srgbClr val="RRGGBB"/> becomes
# 0xBBGGRR here.
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1       000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1       FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2       44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2       E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1   5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2   ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3   A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4   FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5   4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6   70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink     0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink  954F72
